$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 369, pushing the existing rows 369-383 down to 370-384.
$ws.Rows.Item(369).Insert()

# Populate the newly inserted row 369 with the new Kiwi "Primera" price record
# (Terminal Hortofrutícola Agro Chillán, Ñuble), dated 45035 (2023-04-19).
$ws.Range("A369").Value = 7
$ws.Range("B369").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C369").Value = "Ñuble"
$ws.Range("D369").Value = 45035
$ws.Range("E369").Value = 16
$ws.Range("F369").Value = "Fruta"
$ws.Range("G369").Value = 100101
$ws.Range("H369").Value = "Berries"
$ws.Range("I369").Value = 100101007
$ws.Range("J369").Value = "Kiwi"
$ws.Range("K369").Value = "Hayward"
$ws.Range("L369").Value = "Primera"
$ws.Range("M369").Value = 100
$ws.Range("N369").Value = 13000
$ws.Range("O369").Value = 13000
$ws.Range("P369").Value = 13000
$ws.Range("Q369").Value = "`$/bandeja 10 kilos"
$ws.Range("R369").Value = "Región de O'Higgins"
$ws.Range("S369").Value = 1300
$ws.Range("T369").Value = 10
